# "drug usage as virtual column"
# Consolidate the RXQ_RX_J rows in the FILESTREAM sheet: drop the
# "-E" (examination) row entirely, and swap the "-P"/"-T" rows so the
# surviving two rows read "-T" then "-P" (virtual-column re-ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FILESTREAM")

# --- 1. Swap the da-name / data-dict values for rows 3 and 4 -------------
# Row 3 currently holds ...-RXQ_RX_J-P, row 4 holds ...-RXQ_RX_J-T.
# After the edit row 3 should hold -T and row 4 should hold -P.
$valRow3 = $ws.Range("A3").Value2
$valRow4 = $ws.Range("A4").Value2

$ws.Range("A3").Value = $valRow4
$ws.Range("B3").Value = $valRow4
$ws.Range("A4").Value = $valRow3
$ws.Range("B4").Value = $valRow3

# --- 2. Split the merged E3:E5 hyperlink into per-cell hyperlinks --------
# Remove the old hyperlink that spans E3:E5 (it will be replaced with
# individual hyperlinks on E3 and E4 once row 5 disappears).
$ws.Hyperlinks.Item(2).Delete()

# --- 3. Delete row 5 (the "-E" row) entirely, shifting nothing else -----
$ws.Range("A5").EntireRow.Delete()

# --- 4. Re-create the mailto hyperlinks on the two remaining data rows --
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:example@example.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:example@example.com")

# --- 5. Restore the selection the author left behind ---------------------
$ws.Range("C4").Select()
